$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Delete columns from right to left so earlier indices stay valid:
# R (18) - blank spacer before "Cuoc bien"
# Q (17) - "Cuoc noi bo"
# N (14) - "Cuoc o to moi"
# I (9)  - blank
# H (8)  - blank
# D (4)  - blank
# A (1)  - blank spacer at the very start
$ws.Columns.Item(18).Delete()
$ws.Columns.Item(17).Delete()
$ws.Columns.Item(14).Delete()
$ws.Columns.Item(9).Delete()
$ws.Columns.Item(8).Delete()
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(1).Delete()
